# Applies "Schema changed! Fixed #198 / Added Decision Kinds" edit:
#  - Renames the three color rows' ColorName/Options/Hex values to the
#    new "Decision Kind" keys/labels/colors.
#  - Clears the ColorName column (B) for the data rows since colors are
#    no longer referenced by name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Negative decision
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "|CommissionDecisionNegative|"
$ws.Range("D2").Value = "Отрицательное решение"
$ws.Range("E2").Value = "#FF8375"

# Row 3: Neutral decision
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "|CommissionDecisionNeutral|"
$ws.Range("D3").Value = "Нейтральное решение"
$ws.Range("E3").Value = "#FFF293"

# Row 4: Positive decision
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "|CommissionDecisionPositive|"
$ws.Range("D4").Value = "Положительное решение"
$ws.Range("E4").Value = "#4CFF76"
